{"js": "// Replace the date line and every divison-answer cell in the table with the\n// new values below. Mapping is strictly positional (document order) because\n// several old cell values repeat (e.g. \"60\u00f73=20, 0\" occurs twice, and\n// \"23\u00f72=11, 1\" is the new value for two different cells), so matching by old\n// text alone would be ambiguous. Document order of non-empty paragraphs is\n// stable and 1:1 with this list.\nconst replacements = [\n  \"2023-11-27 Monday\",\n  \"46\u00f77=6, 4\",\n  \"36\u00f72=18, 0\",\n  \"54\u00f78=6, 6\",\n  \"71\u00f76=11, 5\",\n  \"15\u00f75=3, 0\",\n  \"73\u00f78=9, 1\",\n  \"94\u00f75=18, 4\",\n  \"38\u00f76=6, 2\",\n  \"66\u00f77=9, 3\",\n  \"96\u00f74=24, 0\",\n  \"50\u00f74=12, 2\",\n  \"50\u00f75=10, 0\",\n  \"23\u00f72=11, 1\",\n  \"77\u00f79=8, 5\",\n  \"23\u00f72=11, 1\",\n  \"69\u00f78=8, 5\",\n  \"57\u00f78=7, 1\",\n  \"26\u00f77=3, 5\",\n  \"82\u00f79=9, 1\",\n  \"21\u00f76=3, 3\",\n  \"53\u00f77=7, 4\",\n  \"76\u00f77=10, 6\",\n  \"96\u00f78=12, 0\",\n  \"42\u00f73=14, 0\",\n  \"74\u00f77=10, 4\",\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet i = 0;\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.text === \"\") {\n    continue;\n  }\n  if (i >= replacements.length) {\n    break;\n  }\n  paragraph.getRange().insertText(replacements[i], \"Replace\");\n  i++;\n}\n\nawait context.sync();\n", "ps1": "# Replace the date line and every division-answer cell in the table with the\n# new values below. Mapping is strictly positional (document order) because\n# several old cell values repeat (e.g. \"60\u00f73=20, 0\" occurs twice, and\n# \"23\u00f72=11, 1\" is the new value for two different cells), so matching by old\n# text alone would be ambiguous. Document order of the non-empty paragraphs\n# (the date paragraph, then each filled-in table cell, left-to-right /\n# top-to-bottom) is stable and 1:1 with this list.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    \"2023-11-27 Monday\",\n    \"46\u00f77=6, 4\",\n    \"36\u00f72=18, 0\",\n    \"54\u00f78=6, 6\",\n    \"71\u00f76=11, 5\",\n    \"15\u00f75=3, 0\",\n    \"73\u00f78=9, 1\",\n    \"94\u00f75=18, 4\",\n    \"38\u00f76=6, 2\",\n    \"66\u00f77=9, 3\",\n    \"96\u00f74=24, 0\",\n    \"50\u00f74=12, 2\",\n    \"50\u00f75=10, 0\",\n    \"23\u00f72=11, 1\",\n    \"77\u00f79=8, 5\",\n    \"23\u00f72=11, 1\",\n    \"69\u00f78=8, 5\",\n    \"57\u00f78=7, 1\",\n    \"26\u00f77=3, 5\",\n    \"82\u00f79=9, 1\",\n    \"21\u00f76=3, 3\",\n    \"53\u00f77=7, 4\",\n    \"76\u00f77=10, 6\",\n    \"96\u00f78=12, 0\",\n    \"42\u00f73=14, 0\",\n    \"74\u00f77=10, 4\"\n)\n\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    # An empty table-cell paragraph's Range.Text is just the cell-mark +\n    # paragraph-mark (\"\\r\\a\"); trim those control characters to find\n    # paragraphs that actually hold visible text.\n    $trimmed = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($trimmed.Length -gt 0) {\n        if ($i -ge $replacements.Length) {\n            break\n        }\n        $p.Range.Text = $replacements[$i]\n        $i++\n    }\n}\n"}
